$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- "声明类" (declaration category) block ----
$ws.Range("A5").Value = "声明类"

$ws.Range("A6").Value = "CON"
$ws.Range("B6").Value = "name"
$ws.Range("C6").Value = "type"
$ws.Range("D6").Value = "value"
$ws.Range("E6").Value = "type包括int和char,最终常量存储在.data中"

$ws.Range("A7").Value = "VAR"
$ws.Range("B7").Value = "name"
$ws.Range("C7").Value = "type"
$ws.Range("D7").Value = "~"
$ws.Range("E7").Value = "type包括int和char,最终变量存储在.data中"

$ws.Range("A8").Value = "ARR"
$ws.Range("B8").Value = "name"
$ws.Range("C8").Value = "type"
$ws.Range("D8").Value = "number"

$ws.Range("A9").Value = "PARA"
$ws.Range("B9").Value = "name"
$ws.Range("C9").Value = "type "
$ws.Range("D9").Value = "~"

$ws.Range("A10").Value = "FUNC"
$ws.Range("B10").Value = "name"
$ws.Range("C10").Value = "type"
$ws.Range("D10").Value = "number"

$ws.Range("A11").Value = "TEMP"
$ws.Range("B11").Value = "name"
$ws.Range("C11").Value = "type"
$ws.Range("D11").Value = "~"

# ---- "运算类" (operation category) block ----
$ws.Range("A13").Value = "运算类"

$ws.Range("A14").Value = "+"
$ws.Range("B14").Value = "op1"
$ws.Range("C14").Value = "op2"
$ws.Range("D14").Value = "result"

$ws.Range("A15").Value = "-"
$ws.Range("B15").Value = "op1"
$ws.Range("C15").Value = "op2"
$ws.Range("D15").Value = "result"

$ws.Range("A16").Value = "*"
$ws.Range("B16").Value = "op1"
$ws.Range("C16").Value = "op2"
$ws.Range("D16").Value = "result"

$ws.Range("A17").Value = "/"
$ws.Range("B17").Value = "op1"
$ws.Range("C17").Value = "op2"
$ws.Range("D17").Value = "result"

# Column A is widened to fit the new labels
$ws.Columns.Item(1).ColumnWidth = 15.14

# Selection moves along with the new data entry
$ws.Range("E14").Select()

# Page setup (paper size / orientation) was touched in this revision
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
